# Chapitre 8 - initialisation : met en gras les titres "Sommaire" et
# "Titre x" (zones de texte "ZoneTexte 5") des diapositives 285 et 286.

$p = $ppt.ActivePresentation

$targetSlideIds = @(285, 286)

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)

    if ($targetSlideIds -notcontains $slide.SlideID) {
        continue
    }

    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)

        if ($shape.Name -eq "ZoneTexte 5" -and $shape.HasTextFrame) {
            $shape.TextFrame.TextRange.Font.Bold = $true
        }
    }
}
